# Replace the numeric month (column C) with its Spanish three-letter
# abbreviation ("Ene.", "Feb.", ... "Dic.") for every data row of Tabla3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_31")

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

# First capture every row's current numeric month (column C), keyed by row.
$rowMonth = @{}
for ($r = 6; $r -le 85; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $rowMonth[$r] = [int]$cell.Value()
}

# Write the text values back out processing month 1 (Ene.) first, then 2, ...
# 12 (Dic.) last, so the new shared-string table is built up in that order
# (Ene., Feb., Mar., ..., Dic.) regardless of row traversal order.
for ($m = 1; $m -le 12; $m++) {
    for ($r = 6; $r -le 85; $r++) {
        if ($rowMonth[$r] -eq $m) {
            $ws.Cells.Item($r, 3).Value = $monthNames[$m]
        }
    }
}
